# Applies the structural relabeling + documentation changes described in the
# commit:
#   1. HEADER sheet: relabel the key/value rows (DOMAIN/CATEGORY become the
#      new row 3/4 labels, SOURCE_ORG/SOURCE_PERSON move down to rows 5/6 and
#      lose their sample values, CATEGORY/SUB_CATEGORY rows are dropped).
#   2. DIVIDEND_ lookup sheet: insert a new "-" sentinel in front of the
#      numeric code list (column A), shifting every existing code down one
#      row so the list keeps its original order but gains a leading blank
#      choice and the trailing 93/GBP row becomes a normal fully populated
#      row instead of split across two rows.
#   3. DIVIDEND sheet's data validation on column C: widen the source range
#      to match the DIVIDEND_ list's new length (14 -> 15 rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. HEADER sheet
# ---------------------------------------------------------------------
$header = $wb.Worksheets.Item("HEADER")

$header.Range("A3").Value = "DOMAIN"
# B3 (GRMFMC000) is unchanged

$header.Range("A4").Value = "CATEGORY"
# B4 (aloumiotis@foo.gr) is unchanged

$header.Range("A5").Value = "SOURCE_ORG"
$header.Range("B5").Value = ""

$header.Range("A6").Value = "SOURCE_PERSON"
$header.Range("B6").Value = ""

# ---------------------------------------------------------------------
# 2. DIVIDEND_ sheet (hidden lookup list) - shift column A down by one row
#    and insert "-" as the new first entry.
# ---------------------------------------------------------------------
$dividendLookup = $wb.Worksheets.Item("DIVIDEND_")

$codes = @("10", "20", "11", "12", "13", "21", "22", "23", "31", "32", "33", "91", "92", "93")

# Write from the bottom up so we never overwrite a value we still need to
# read (not strictly required since we're using a literal list, but keeps
# the intent obvious and mirrors how the shift actually happens).
for ($i = $codes.Length - 1; $i -ge 0; $i--) {
    $dividendLookup.Cells.Item($i + 2, 1).Value = $codes[$i]
}
$dividendLookup.Cells.Item(1, 1).Value = "-"

# ---------------------------------------------------------------------
# 3. DIVIDEND sheet - widen the C-column validation list to the new range
# ---------------------------------------------------------------------
$dividend = $wb.Worksheets.Item("DIVIDEND")
$dividend.Range("C4:C20").Validation.Formula1 = "'DIVIDEND_'!`$A`$1:`$A`$15"
